$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "abc5@gmail.com"
$ws.Range("E3").Value = "abc6@gmail.com"
$ws.Range("E4").Value = "abc7@gmail.com"
$ws.Range("E5").Value = "abc8@gmail.com"
